# quarterly.xlsx update — shift the quarterly columns (drop the oldest
# quarter "فصل دوم منتهی به 1399/06" and append the newest quarter
# "فصل چهارم منتهی به 1401/12"), and refresh the figures for every
# quarterly metric row accordingly (commit: "update database and change
# read_price algorithm").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param(
        [int]$Row,
        [object[]]$Values
    )
    $count = $Values.Length
    $arr = New-Object 'object[,]' 1, $count
    for ($i = 0; $i -lt $count; $i++) {
        $arr[0, $i] = $Values[$i]
    }
    $startCell = $ws.Cells.Item($Row, 5)   # column E
    $endCell   = $ws.Cells.Item($Row, 4 + $count)
    $rng = $ws.Range($startCell, $endCell)
    $rng.Value = $arr
}

# --- Header rows (quarter labels), columns E:N -----------------------------
$quarterLabels = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
Set-RowValues 8 $quarterLabels
Set-RowValues 24 $quarterLabels

# --- Data rows, columns E:N -------------------------------------------------
Set-RowValues 10 @(129143, 185963, 88712, 198575, 184666, 324616, 188412, 313155, 310412, 422682)
Set-RowValues 13 @(102859, 260183, 64638, 284744, 122262, 312886, 121780, 255330, 369475, 488076)
Set-RowValues 14 @(15989, -4056, 2607, 11546, 13775, 23894, 8906, 4987, 3584, 4121)
Set-RowValues 15 @(1160, 463, 643, 724, 821, 1010, 844, 1607, 1577, 1270)
Set-RowValues 16 @(2380, 3953, 5258, 5400, 6063, 6548, 6990, 6866, 4870, 6618)
Set-RowValues 17 @(70886, 38760, 76812, 59782, 89512, 100278, 155205, 106213, 148729, 134777)
Set-RowValues 19 @(58696, 32052, 37982, 102696, 57832, 76464, 232159, -52357, 421954, -332862)
Set-RowValues 20 @(381113, 517318, 276652, 663467, 474931, 845696, 714296, 635801, 1260601, 724682)
Set-RowValues 26 @(224, 212, 214, 219, 219, 741, 211, 240, 220, 741)
Set-RowValues 27 @(512, 550, 640, 617, 617, 240, 735, 741, 803, 240)
